$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C slightly (6.25 -> 6.5 as stored in the sheet XML).
# The COM layer adds its own fixed padding/offset when writing the <col>
# width out to XML, so the value set here must be pre-compensated.
$ws.Columns.Item(3).ColumnWidth = 5.666666666666667

# Fill in the missing X7/Y7 values for the existing last row
$ws.Range("X7").Value = 0.42000000000000171
$ws.Range("Y7").Value = "Up"

# Append new row 8 of data
$ws.Range("A8").Value = 42649.879930555559
$ws.Range("A8").NumberFormat = "m/d/yy h:mm"
$ws.Range("B8").Value = -3
$ws.Range("C8").Value = "Neutral"
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = 4245
$ws.Range("F8").Value = 235
$ws.Range("G8").Value = 43
$ws.Range("H8").Value = 55
$ws.Range("I8").Value = 69
$ws.Range("J8").Value = 29
$ws.Range("K8").Value = 7776
$ws.Range("L8").Value = 40
$ws.Range("M8").Value = 51
$ws.Range("N8").Value = 21
$ws.Range("O8").Value = 9
$ws.Range("P8").Value = "Named"
$ws.Range("Q8").Value = 47.321424984051369
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0.0834
$ws.Range("S8").NumberFormat = "0.00%"
$ws.Range("T8").Value = -0.0062
$ws.Range("T8").NumberFormat = "0.00%"
$ws.Range("U8").Value = 2.31
$ws.Range("V8").Value = "N/A"
$ws.Range("W8").Value = 0
